$wb = $excel.ActiveWorkbook
$pwm = $wb.Worksheets.Item("PWM")
$radio = $wb.Worksheets.Item("Radio")

# --- Header row for the new "Colors" block (mirrors the "Experimental" header on Radio sheet) ---
$pwm.Range("A12").Value = "Colors"
$pwm.Range("A12").Style = "Accent1"

# --- RED/GREEN/BLUE PWM max + color step (mirrors the "Good" style rows) ---
$pwm.Range("A13").Value = "RED_PWM_MAX"
$pwm.Range("B13").Value = 160
$pwm.Range("A13:B13").Style = "Good"

$pwm.Range("A14").Value = "GREEN_PWM_MAX"
$pwm.Range("B14").Value = 160
$pwm.Range("A14:B14").Style = "Good"

$pwm.Range("A15").Value = "BLUE_PWM_MAX"
$pwm.Range("B15").Value = 160
$pwm.Range("A15:B15").Style = "Good"

$pwm.Range("A16").Value = "COLOR_STEP"
$pwm.Range("B16").Value = 8
$pwm.Range("A16:B16").Style = "Good"

# --- Variations per color (mirrors the "Neutral" style) ---
$pwm.Range("A17").Value = "Red variations"
$pwm.Range("B17").Formula = '=B13/$B$16'
$pwm.Range("A17:B17").Style = "Neutral"

$pwm.Range("A18").Value = "Green variations"
$pwm.Range("B18").Formula = '=B14/$B$16'
$pwm.Range("A18:B18").Style = "Neutral"

$pwm.Range("A19").Value = "Blue variations"
$pwm.Range("B19").Formula = '=B15/$B$16'
$pwm.Range("A19:B19").Style = "Neutral"

# --- Total number of colors (mirrors the "Calculation" style) ---
$pwm.Range("A20").Value = "Colors count"
$pwm.Range("B20").Formula = '=B19*B18*B17'
$pwm.Range("A20:B20").Style = "Calculation"

# --- View / selection changes ---
$pwm.Range("B21").Select()
$radio.Range("B12").Select()
$pwm.Activate()

$wb.Save()
